$d = $word.ActiveDocument

# The "Requisitos" bulleted list paragraph currently holds three runs, each
# "<text><w:br/>":
#   1) LOM3022 -  Materiais para a  Indústria Química  (Requisito fraco)
#   2) LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)
#   3) LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito fraco)
#
# The commit moves the LOB1009 requirement to the front, leaving the other
# two requirements in their original relative order:
#   1) LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito fraco)
#   2) LOM3022 -  Materiais para a  Indústria Química  (Requisito fraco)
#   3) LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)

$textLOM = "LOM3022 -  Materiais para a  Indústria Química  (Requisito fraco)"
$textLOB = "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito fraco)"

# Chr(11) is Word's manual line break character (serializes as <w:br/>),
# the same element already used to separate these requirement lines.
$lineBreak = [char]11

# 1) Insert a brand new "LOB1009...<break>" run right before the LOM3022
#    run. Using InsertBefore on a collapsed range creates an independent run
#    instead of merging with its neighbour.
$lomRange = $d.Content
$lomRange.Find.Execute($textLOM, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$insertionPoint = $d.Range($lomRange.Start, $lomRange.Start)
$insertionPoint.InsertBefore($textLOB + $lineBreak)

# 2) Remove the original (now duplicated) "LOB1009...<break>" run that used
#    to sit at the end of the list, by locating the occurrence that comes
#    after the one we just inserted.
$firstOccurrence = $d.Content
$firstOccurrence.Find.Execute($textLOB, $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null

$remainder = $d.Range($firstOccurrence.End, $d.Content.End)
$remainder.Find.Execute($textLOB, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$oldRun = $d.Range($remainder.Start, $remainder.End + 1)
$oldRun.Delete()
